# "update scripts wuth new tpm" -- the NATMI L1cam/Alcam ligand-receptor
# table was recomputed from updated TPM source data. The sending/target
# cluster grid used to be the full 3x3 (ECs/FAPs/MuSCs x ECs/FAPs/MuSCs,
# rows 2-10); the refreshed run drops every row whose Target cluster is
# "ECs" (old rows 2, 5 and 8) and recalculates every numeric column for
# the remaining 6 sending/target combinations (now rows 2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three rows for Target cluster "ECs" - everything below shifts
# up and the sheet's used range shrinks from A1:T10 to A1:T7.
$ws.Rows("8:10").Delete()

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 7.741029
$ws.Range("H2").Value = 23.223087
$ws.Range("I2").Value = 0.4930486933812723
$ws.Range("J2").Value = 0.4930486933812723
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.4771596666666666
$ws.Range("N2").Value = 1.431479
$ws.Range("O2").Value = 0.4973652976730675
$ws.Range("P2").Value = 0.4973652976730676
$ws.Range("Q2").Value = 3.693706817297
$ws.Range("R2").Value = 33.243361355673
$ws.Range("S2").Value = 0.2452253101508935
$ws.Range("T2").Value = 0.2452253101508935

# Row 3: ECs -> MuSCs
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 7.741029
$ws.Range("H3").Value = 23.223087
$ws.Range("I3").Value = 0.4930486933812723
$ws.Range("J3").Value = 0.4930486933812723
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.482215
$ws.Range("N3").Value = 1.446645
$ws.Range("O3").Value = 0.5026347023269324
$ws.Range("P3").Value = 0.5026347023269324
$ws.Range("Q3").Value = 3.732840299235
$ws.Range("R3").Value = 33.595562693115
$ws.Range("S3").Value = 0.2478233832303788
$ws.Range("T3").Value = 0.2478233832303788

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 0.3045986666666667
$ws.Range("H4").Value = 0.9137960000000001
$ws.Range("I4").Value = 0.0194007766416684
$ws.Range("J4").Value = 0.0194007766416684
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.4771596666666666
$ws.Range("N4").Value = 1.431479
$ws.Range("O4").Value = 0.4973652976730675
$ws.Range("P4").Value = 0.4973652976730676
$ws.Range("Q4").Value = 0.1453421982537778
$ws.Range("R4").Value = 1.308079784284
$ws.Range("S4").Value = 0.0096492730494721
$ws.Range("T4").Value = 0.009649273049472102

# Row 5: FAPs -> MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.3045986666666667
$ws.Range("H5").Value = 0.9137960000000001
$ws.Range("I5").Value = 0.0194007766416684
$ws.Range("J5").Value = 0.0194007766416684
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.482215
$ws.Range("N5").Value = 1.446645
$ws.Range("O5").Value = 0.5026347023269324
$ws.Range("P5").Value = 0.5026347023269324
$ws.Range("Q5").Value = 0.1468820460466667
$ws.Range("R5").Value = 1.32193841442
$ws.Range("S5").Value = 0.009751503592196301
$ws.Range("T5").Value = 0.009751503592196301

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 7.654706000000001
$ws.Range("H6").Value = 22.964118
$ws.Range("I6").Value = 0.4875505299770593
$ws.Range("J6").Value = 0.4875505299770593
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4771596666666666
$ws.Range("N6").Value = 1.431479
$ws.Range("O6").Value = 0.4973652976730675
$ws.Range("P6").Value = 0.4973652976730676
$ws.Range("Q6").Value = 3.652516963391334
$ws.Range("R6").Value = 32.872652670522
$ws.Range("S6").Value = 0.242490714472702
$ws.Range("T6").Value = 0.242490714472702

# Row 7: MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 7.654706000000001
$ws.Range("H7").Value = 22.964118
$ws.Range("I7").Value = 0.4875505299770593
$ws.Range("J7").Value = 0.4875505299770593
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.482215
$ws.Range("N7").Value = 1.446645
$ws.Range("O7").Value = 0.5026347023269324
$ws.Range("P7").Value = 0.5026347023269324
$ws.Range("Q7").Value = 3.69121405379
$ws.Range("R7").Value = 33.22092648411
$ws.Range("S7").Value = 0.2450598155043573
$ws.Range("T7").Value = 0.2450598155043573
